$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")

# Insert 4 new columns before column D (old D/E/F shift to H/I/J),
# carrying forward the formatting that was in column D.
$ws.Range("D1:G1").EntireColumn.Insert()

# Set header values for the newly inserted columns.
$ws.Range("D1").Value = "CUST.REMARKS:1"
$ws.Range("E1").Value = "INTEND.DATE"
$ws.Range("F1").Value = "EXP.DATE"
$ws.Range("G1").Value = "TAX.INTEREST.TYPE:1"

# Reset the selection to match the saved view state.
$ws.Range("I8").Select()
